$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 227 (old row 227 "Tercera/44544" and
# old row 228 "Primera/44160" shift down to 229 and 230 respectively).
$ws.Rows.Item(227).Insert()
$ws.Rows.Item(227).Insert()

# New row 227: Primera, 03-02-2022 (44595), 400 vol, 3000 min/max/avg, $/unidad, O'Higgins
$ws.Cells.Item(227, 1).Value = 4
$ws.Cells.Item(227, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(227, 3).Value = "Los Lagos"
$ws.Cells.Item(227, 4).Value = 44595
$ws.Cells.Item(227, 5).Value = 10
$ws.Cells.Item(227, 6).Value = 100112028
$ws.Cells.Item(227, 7).Value = "Sandia"
$ws.Cells.Item(227, 8).Value = "Sin especificar"
$ws.Cells.Item(227, 9).Value = "Primera"
$ws.Cells.Item(227, 10).Value = 400
$ws.Cells.Item(227, 11).Value = 3000
$ws.Cells.Item(227, 12).Value = 3000
$ws.Cells.Item(227, 13).Value = 3000
$ws.Cells.Item(227, 14).Value = "`$/unidad"
$ws.Cells.Item(227, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(227, 16).Value = 3000
$ws.Cells.Item(227, 17).Value = 1
$ws.Cells.Item(227, 18).Value = "Hortaliza"

# New row 228: Segunda, 03-02-2022 (44595), 400 vol, 2500 min/max/avg, $/unidad, O'Higgins
$ws.Cells.Item(228, 1).Value = 4
$ws.Cells.Item(228, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(228, 3).Value = "Los Lagos"
$ws.Cells.Item(228, 4).Value = 44595
$ws.Cells.Item(228, 5).Value = 10
$ws.Cells.Item(228, 6).Value = 100112028
$ws.Cells.Item(228, 7).Value = "Sandia"
$ws.Cells.Item(228, 8).Value = "Sin especificar"
$ws.Cells.Item(228, 9).Value = "Segunda"
$ws.Cells.Item(228, 10).Value = 400
$ws.Cells.Item(228, 11).Value = 2500
$ws.Cells.Item(228, 12).Value = 2500
$ws.Cells.Item(228, 13).Value = 2500
$ws.Cells.Item(228, 14).Value = "`$/unidad"
$ws.Cells.Item(228, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(228, 16).Value = 2500
$ws.Cells.Item(228, 17).Value = 1
$ws.Cells.Item(228, 18).Value = "Hortaliza"

Write-Output "done"
